$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-26 10:07:06"
$wsZhCn.Range("G5").Value = "2016-01-26 10:07:49"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-26 10:07:17"
$wsDeDe.Range("G5").Value = "2016-01-26 10:08:10"
